# Update "想去人数" (column F) counts on the "展览" and "全部类型" sheets
# to reflect newly scraped totals.

$wb = $excel.ActiveWorkbook

# Row -> new F value, per worksheet name.
$updates = @{
    "展览" = @{
        2  = 4704
        3  = 2554
        4  = 71
        6  = 243
        7  = 142
        8  = 220
        9  = 185
        10 = 1811
        11 = 327
        12 = 4165
        13 = 51
    }
    "全部类型" = @{
        2  = 4704
        3  = 2554
        5  = 71
        8  = 243
        9  = 142
        10 = 220
        11 = 185
        14 = 1811
        15 = 327
        16 = 4165
        17 = 51
    }
}

foreach ($sheetName in $updates.Keys) {
    $ws = $wb.Worksheets.Item($sheetName)
    $rowsForSheet = $updates[$sheetName]
    foreach ($row in $rowsForSheet.Keys) {
        $ws.Cells.Item($row, 6).Value = $rowsForSheet[$row]
    }
}
